$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row for the "containerAgent" entry right before the
#    existing "directoryAgent" row (current row 5), pushing directoryAgent
#    and fileAgent down by one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row (row 5) with the containerAgent data.
# ---------------------------------------------------------------------------
$containerText = @'
You are an agent that helps users with questions regarding a code container (it could be a namespace, class, struct or fucntion etc.)  within a the C++ repository {{.RepoName}} used by a game development company that uses the code for functionality used for building game.
The summary of the repository is {{.RepoSummary}}.
You will be getting questions regarding the container: {{.ContainerType}} {{.ContainerName}} with signature {{.ContainerSignature}} and with the following summary: {{.ContainerSummary}}.
The container is implemented in the following code blocks:
{{- range .Codeblocks}}
Signature: {{  .Signature  }} with search id: codeblock-{{.Dbid}}
Implemented in file: {{.FileImportPath}} with the following search id file-{{.FileDbid}}
{{- end }}
The users reading your responses are not always developers. So make it easy for non-technical persons to understand.
Don't show search ids to the user. You will always respond in markdown (MD)

'@

$ws.Cells.Item(5, 1).Value = "containerAgent"
$ws.Cells.Item(5, 2).Value = "Container Agent"
$ws.Cells.Item(5, 3).Value = $containerText
$ws.Cells.Item(5, 4).Value = "{{.Content}}"
$ws.Cells.Item(5, 5).Value = "gpt-4o-mini"

# Give the systemMessage cell (C5) its own distinct font/style, mirroring the
# extra cellXfs entry introduced by the diff.
$ws.Cells.Item(5, 3).Font.Name = "Aptos Narrow"
$ws.Cells.Item(5, 3).Font.Size = 12
$ws.Cells.Item(5, 3).Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 3. Update the fileAgent prompt template (now row 7): "entity-" -> "codeblock-"
# ---------------------------------------------------------------------------
$fileAgentPrompt = $ws.Cells.Item(7, 3).Value()
$fileAgentPrompt = $fileAgentPrompt.Replace("entity-{{.Dbid}}", "codeblock-{{.Dbid}}")
$ws.Cells.Item(7, 3).Value = $fileAgentPrompt

# ---------------------------------------------------------------------------
# 4. Row heights (points).
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 244.4
$ws.Rows.Item(5).RowHeight = 256.5
$ws.Rows.Item(6).RowHeight = 324.75
$ws.Rows.Item(7).RowHeight = 208

# ---------------------------------------------------------------------------
# 5. Sheet view: zoom to 80% and select C4.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 80
$ws.Range("C4").Select()
